# Remove the paragraph that follows the "MRK" heading and contains the
# lone italic "Marcos" run (the paragraph's text and its paragraph mark
# are both deleted, merging it away so the document structure matches
# the post-edit OOXML).
$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Marcos" -and $p.Style.NameLocal -eq "Normal") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
